$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G: rename "Trizol" -> "TRIzol" and drop the now-unneeded
#     style/validation (cells revert to the default "Normal" look) ---
$ws.Range("G2:G27").Style = "Normal"
$ws.Range("G2:G27").Value = "TRIzol"
$ws.Range("G2:G27").Validation.Delete()

# --- Column H: convert the yes/no text answers into a real boolean,
#     formatted to still read as TRUE/FALSE, and drop the old validation ---
$ws.Range("H2:H27").Style = "Normal"
$ws.Range("H2:H27").Value = $false
$ws.Range("H2:H27").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("H2:H27").Validation.Delete()

# --- Row heights for the data rows shrink from 16 to 15 ---
$ws.Range("2:27").RowHeight = 15

# --- Selection moves to the column that was just edited ---
$ws.Range("G2:G27").Select()
